$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7515
$ws.Range("J69").Value = 7515
$ws.Range("L69").Value = 22545
$ws.Range("N69").Value = -24293
$ws.Range("H72").Value = 7515
$ws.Range("J72").Value = 7515
$ws.Range("L72").Value = 67635
$ws.Range("N72").Value = -76371
$ws.Range("H92").Value = 1624.75
$ws.Range("J92").Value = 997.5
$ws.Range("L92").Value = 997.5
$ws.Range("N92").Value = -3493.5
$ws.Range("H100").Value = 2800.1333
$ws.Range("I100").Value = 2227
$ws.Range("J100").Value = 4376.25
$ws.Range("K100").Value = 2227
$ws.Range("L100").Value = 4376.25
$ws.Range("M100").Value = -1686
$ws.Range("N100").Value = -5458.25
$ws.Range("H112").Value = 4032.88
$ws.Range("J112").Value = 4038.4167
$ws.Range("L112").Value = 12115.2501
$ws.Range("N112").Value = -14331.2501
$ws.Range("H137").Value = 23815360
$ws.Range("I137").Value = 83335000
$ws.Range("J137").Value = 7505.2
$ws.Range("K137").Value = 250005000
$ws.Range("L137").Value = 22515.6
$ws.Range("M137").Value = -250002450
$ws.Range("N137").Value = -27615.6
$ws.Range("H138").Value = 2689.0588
$ws.Range("J138").Value = 3893.5
$ws.Range("L138").Value = 11680.5
$ws.Range("N138").Value = -21960.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1769.091
$ws.Range("I45").Value = 1714
$ws.Range("J45").Value = 1815
$ws.Range("K45").Value = 1714
$ws.Range("L45").Value = 1815
$ws.Range("M45").Value = -1337
$ws.Range("N45").Value = -2569
$ws.Range("H61").Value = 4768438.5
$ws.Range("I61").Value = 8003.154
$ws.Range("J61").Value = 12504146
$ws.Range("K61").Value = 8003.154
$ws.Range("L61").Value = 12504146
$ws.Range("M61").Value = -7791.154
$ws.Range("N61").Value = -12504570
$ws.Range("H122").Value = 1549.25
$ws.Range("I122").Value = 1349
$ws.Range("K122").Value = 4047
$ws.Range("M122").Value = -1597
$ws.Range("H132").Value = 3049.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3049.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9148.5
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 4768438.5
$ws.Range("I136").Value = 8003.154
$ws.Range("J136").Value = 12504146
$ws.Range("K136").Value = 24009.462
$ws.Range("L136").Value = 37512438
$ws.Range("M136").Value = -21459.462
$ws.Range("N136").Value = -37517538

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 167.5
$ws.Range("I22").Value = 167.5
$ws.Range("K22").Value = 167.5
$ws.Range("M22").Value = 5.5
$ws.Range("H94").Value = 3231.35
$ws.Range("J94").Value = 2033
$ws.Range("L94").Value = 2033
$ws.Range("N94").Value = -2935

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2928836.5
$ws.Range("I31").Value = 3707193
$ws.Range("K31").Value = 3707193
$ws.Range("M31").Value = -3706898
$ws.Range("H34").Value = 2928836.5
$ws.Range("I34").Value = 3707193
$ws.Range("K34").Value = 3707193
$ws.Range("M34").Value = -3706991
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H58").Value = 4640612.5
$ws.Range("J58").Value = 5567519.5
$ws.Range("L58").Value = 5567519.5
$ws.Range("N58").Value = -5567925.5
$ws.Range("H62").Value = 4191
$ws.Range("J62").Value = 4229.2
$ws.Range("L62").Value = 4229.2
$ws.Range("N62").Value = -5477.2
$ws.Range("H65").Value = 4191
$ws.Range("J65").Value = 4229.2
$ws.Range("L65").Value = 21146
$ws.Range("N65").Value = -27386
$ws.Range("H70").Value = 39125
$ws.Range("J70").Value = 39125
$ws.Range("L70").Value = 39125
$ws.Range("N70").Value = -39755
$ws.Range("H73").Value = 39125
$ws.Range("J73").Value = 39125
$ws.Range("L73").Value = 39125
$ws.Range("N73").Value = -41309
$ws.Range("H93").Value = 12865.9
$ws.Range("I93").Value = 12865.9
$ws.Range("K93").Value = 12865.9
$ws.Range("M93").Value = -10993.9
$ws.Range("H103").Value = 10295.1
$ws.Range("I103").Value = 7883.4443
$ws.Range("K103").Value = 7883.4443
$ws.Range("M103").Value = -6711.4443
$ws.Range("H105").Value = 8762.643
$ws.Range("I105").Value = 9882.362999999999
$ws.Range("K105").Value = 9882.362999999999
$ws.Range("M105").Value = -8135.362999999999
$ws.Range("H136").Value = 4640612.5
$ws.Range("J136").Value = 5567519.5
$ws.Range("L136").Value = 16702558.5
$ws.Range("N136").Value = -16707658.5
$ws.Range("H141").Value = 292456.1
$ws.Range("J141").Value = 335224.22
$ws.Range("L141").Value = 335224.22
$ws.Range("N141").Value = -345584.22

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 139.77777
$ws.Range("I2").Value = 85.07692
$ws.Range("J2").Value = 282
$ws.Range("K2").Value = 510.46152
$ws.Range("L2").Value = 1692
$ws.Range("M2").Value = -397.46152
$ws.Range("N2").Value = -1918
$ws.Range("H64").Value = 11749.5
$ws.Range("I64").Value = 6498
$ws.Range("J64").Value = 13500
$ws.Range("K64").Value = 19494
$ws.Range("L64").Value = 40500
$ws.Range("M64").Value = -19224
$ws.Range("N64").Value = -41040
$ws.Range("H67").Value = 11749.5
$ws.Range("I67").Value = 6498
$ws.Range("J67").Value = 13500
$ws.Range("K67").Value = 19494
$ws.Range("L67").Value = 40500
$ws.Range("M67").Value = -18558
$ws.Range("N67").Value = -42372
$ws.Range("H93").Value = 3000
$ws.Range("I93").Value = 3000
$ws.Range("K93").Value = 9000
$ws.Range("M93").Value = -7128

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 12345
$ws.Range("J4").Value = 12345
$ws.Range("L4").Value = 12345
$ws.Range("N4").Value = -12569
$ws.Range("H107").Value = 8023
$ws.Range("I107").Value = 10147.546
$ws.Range("K107").Value = 10147.546
$ws.Range("M107").Value = -8227.546
$ws.Range("H126").Value = 25417.334
$ws.Range("I126").Value = 25417.334
$ws.Range("K126").Value = 76252.00199999999
$ws.Range("M126").Value = -73782.00199999999
$ws.Range("H132").Value = 16373.435
$ws.Range("J132").Value = 100000
$ws.Range("L132").Value = 300000
$ws.Range("N132").Value = -305060

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2408.6667
$ws.Range("I22").Value = 2337.6
$ws.Range("K22").Value = 2337.6
$ws.Range("M22").Value = -2042.6
$ws.Range("H27").Value = 2408.6667
$ws.Range("I27").Value = 2337.6
$ws.Range("K27").Value = 2337.6
$ws.Range("M27").Value = -2230.6
$ws.Range("H55").Value = 819.9545000000001
$ws.Range("I55").Value = 832.93335
$ws.Range("K55").Value = 832.93335
$ws.Range("M55").Value = -659.93335

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1999
$ws.Range("J3").Value = 1999
$ws.Range("L3").Value = 1999
$ws.Range("N3").Value = -2227
$ws.Range("H51").Value = 19975.625
$ws.Range("I51").Value = 15114.8
$ws.Range("J51").Value = 28077
$ws.Range("K51").Value = 15114.8
$ws.Range("L51").Value = 28077
$ws.Range("M51").Value = -14604.8
$ws.Range("N51").Value = -29097
$ws.Range("H100").Value = 3693
$ws.Range("I100").Value = 2775.1428
$ws.Range("K100").Value = 5550.2856
$ws.Range("M100").Value = -5009.2856
$ws.Range("H113").Value = 1530.3572
$ws.Range("J113").Value = 1610.1818
$ws.Range("L113").Value = 4830.5454
$ws.Range("N113").Value = -9170.545399999999
$ws.Range("H122").Value = 70373.12
$ws.Range("I122").Value = 4470.4375
$ws.Range("K122").Value = 13411.3125
$ws.Range("M122").Value = -10961.3125
$ws.Range("H126").Value = 2621.3809
$ws.Range("I126").Value = 2680.3125
$ws.Range("K126").Value = 8040.9375
$ws.Range("M126").Value = -5570.9375
$ws.Range("H132").Value = 3789580.8
$ws.Range("I132").Value = 4168095
$ws.Range("K132").Value = 12504285
$ws.Range("M132").Value = -12501755
